$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Who is the president of the USA?"
$ws.Range("B6").Value = "Joe Biden"
$ws.Range("C6").Value = "Person"

$ws.Range("A5").Value = "When did Miachel Schumacher win his first F1 World Drivers Title?"
$ws.Range("B5").Value = 1994
$ws.Range("C5").Value = "Year"

$ws.Range("B7").Value = "Max Verstappen"
$ws.Range("A7").Value = "Who was the F1 World Champion in 2022?"
$ws.Range("C7").Value = "Person"

$ws.Columns.Item(2).EntireColumn.AutoFit()

[void]$ws.Range("D14").Select()
